$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "currentPhase"
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
